$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("4/4/2019", "654", "654", "654", "654", "654", "654", "654", "654", "654"),
    @("4/4/2019", "6987", "9456", "453", "15849", "84", "489", "4968", "484", "44"),
    @("4/4/2019", "5", "984", "65156", "484", "415", "89", "4", "51", "684894"),
    @("4/4/2019", "654", "65", "878", "7979", "77", "98", "79", "8798", "7987"),
    @("4/4/2019", "6", "897", "987", "8", "79", "97", "7798", "798", "7987"),
    @("4/4/2019", "654", "564", "6456", "564", "54", "56465", "456", "4654", "54")
)

$startRow = 35
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = "'" + $rowData[$j]
    }
}

$ws.Range("A35:J40").Style = "Normal"
